$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "Item Proficiency Requirements" table: widen it, add a 4th column
#    ("Max Save DC") with per-row Save DC values.
# ---------------------------------------------------------------------------
$t = $d.Tables.Item(1)

# Add the new 4th column, then size every column to the target dxa widths
# (values below are in points = dxa / 20).
$t.Columns.Add() | Out-Null
$t.Columns.Item(1).Width = 118.15   # 2363 dxa
$t.Columns.Item(2).Width = 148.5    # 2970 dxa
$t.Columns.Item(3).Width = 180.0    # 3600 dxa
$t.Columns.Item(4).Width = 99.0     # 1980 dxa

# Overall table width + (negative) left indent
$t.PreferredWidth = 545.65          # 10913 dxa
$t.Rows.LeftIndent = -39.25         # -785 dxa

# Fill in the new column's cells.
$maxSaveValues = @("Max Save DC", "13", "15", "17", "18", "19+")
for ($row = 1; $row -le 6; $row++) {
    $cell = $t.Rows.Item($row).Cells.Item(4)
    $cell.Range.Text = $maxSaveValues[$row - 1]

    # Re-fetch after the structural edit above before touching formatting.
    $cell = $t.Rows.Item($row).Cells.Item(4)
    $pr = $cell.Range.Paragraphs.Item(1).Range
    $pr.Font.Size = 12
    $pr.Font.SizeBi = 12
    if ($row -eq 1) {
        $pr.Font.Bold = 1
    }
    $pr.ParagraphFormat.Alignment = 1
}

# ---------------------------------------------------------------------------
# 2. "Blacksmithing Notes" heading: merge the two runs ("Blacksmithing" +
#    " Notes") into a single run, dropping the lastRenderedPageBreak.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Blacksmithing Notes", $true, $false, $false, $false, `
    $false, $true, 1, $false, "Blacksmithing Notes", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Whetstones line: reword to "Whetstones: (Cost: 100 x Level gp)" with the
#    cost formula moved earlier, drop the stray _GoBack bookmark, and keep
#    "gp" / ")" as separate trailing runs.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Whetstones: (Cost: 100 x ", $true, $false, $false, `
    $false, $false, $true, 1, $false, "Whetstones: (Cost: 100 x Level ", 2) | Out-Null
$d.Content.Find.Execute("Level gp)", $true, $false, $false, $false, $false, `
    $true, 1, $false, "gp)", 2) | Out-Null

$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# Split the trailing "gp)" run into "gp" and ")" runs (matching the target
# markup, which wraps "gp" in proofErr spell-check tags).
$f = $d.Content
$f.Find.Execute("gp)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$closeParen = $d.Range($f.Start + 2, $f.End)
$closeParen.Font.Bold = 0
$closeParen.Font.Bold = 1

Write-Host "Edit complete"
